$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p017r_1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p017r_1</id>", 2)
